$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = '29.429.57'
$ws.Range("E2").Value = '  +0.36%  '
$ws.Range("D3").Value = '1.869.04'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '243.89'
$ws.Range("E5").Value = '  +0.73%  '
$ws.Range("D6").Value = '0.7069'
$ws.Range("E6").Value = '  -0.53%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '0.07890'
$ws.Range("E8").Value = '  -1.29%  '
$ws.Range("D9").Value = '0.3134'
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("D10").Value = '24.57'
$ws.Range("E10").Value = '  -1.52%  '
$ws.Range("E11").Value = '  -4.06%  '
$ws.Range("D12").Value = '1.899.58'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '5.202'
$ws.Range("E13").Value = '  -0.77%  '
$ws.Range("D14").Value = '93.46'
$ws.Range("E14").Value = '  -1.02%  '
$ws.Range("D15").Value = '0.7020'
$ws.Range("E15").Value = '  -1.26%  '
$ws.Range("D16").Value = '6.515'
$ws.Range("E16").Value = '  +2.47%  '
$ws.Range("D17").Value = '29.688.70'
$ws.Range("E17").Value = '  +1.18%  '
$ws.Range("D18").Value = '0.000008361'
$ws.Range("E18").Value = '  -1.68%  '
$ws.Range("D19").Value = '252.21'
$ws.Range("E19").Value = '  +2.85%  '
$ws.Range("D20").Value = '2.157.26'
$ws.Range("E20").Value = '  +0.70%  '
$ws.Range("D21").Value = '13.11'
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("D22").Value = '1.002'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '7.662'
$ws.Range("E23").Value = '  -1.49%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '0.1555'
$ws.Range("E25").Value = '  +0.14%  '
$ws.Range("D26").Value = '8.998'
$ws.Range("E26").Value = '  -0.56%  '
$ws.Range("D27").Value = '161.57'
$ws.Range("E27").Value = '  -0.59%  '
$ws.Range("E28").Value = '  +1.03%  '
$ws.Range("E29").Value = '  -0.24%  '
$ws.Range("D30").Value = '4.330'
$ws.Range("E30").Value = '  -1.91%  '
$ws.Range("D31").Value = '4.255'
$ws.Range("E31").Value = '  -1.34%  '
$ws.Range("D32").Value = '1.204'
$ws.Range("E32").Value = '  +1.73%  '
$ws.Range("D33").Value = '0.05318'
$ws.Range("E33").Value = '  -0.96%  '
$ws.Range("D34").Value = '1.896'
$ws.Range("E34").Value = '  -1.92%  '
$ws.Range("D35").Value = '0.7478'
$ws.Range("E35").Value = '  -2.17%  '
$ws.Range("D36").Value = '1.172'
$ws.Range("E36").Value = '  -0.93%  '
$ws.Range("D37").Value = '2.713'
$ws.Range("E37").Value = '  +0.92%  '
$ws.Range("D38").Value = '0.01891'
$ws.Range("E38").Value = '  +0.60%  '
$ws.Range("D39").Value = '1.278.29'
$ws.Range("E39").Value = '  +1.71%  '
$ws.Range("D40").Value = '2.769'
$ws.Range("E40").Value = '  +0.66%  '
$ws.Range("D41").Value = '0.8951'
$ws.Range("E41").Value = '  -2.23%  '
$ws.Range("D42").Value = '6.091'
$ws.Range("E42").Value = '  -6.46%  '
$ws.Range("D43").Value = '109.33'
$ws.Range("E43").Value = '  -3.07%  '
$ws.Range("D44").Value = '71.38'
$ws.Range("E44").Value = '  -3.77%  '
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  -0.03%  '
$ws.Range("D46").Value = '0.00000000129'
$ws.Range("E46").Value = '  -2.65%  '
$ws.Range("D47").Value = '2.033.02'
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").Value = '1.796'
$ws.Range("E48").Value = '  -0.12%  '
$ws.Range("B49").Value = 'Mantle'
$ws.Range("C49").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D49").Value = '0.5182'
$ws.Range("E49").Value = '  -0.77%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '9.542'
$ws.Range("E50").Value = '  +0.98%  '
$ws.Range("D51").Value = '0.4311'
$ws.Range("E51").Value = '  -1.26%  '
